# Update the lesson catalog on Sheet1:
#  - drop the old "Con" (tag list) column entirely
#  - simplify the Quiz/Worksheet/Flashcard link columns to hold direct links
#  - add a new Grade 7 / Science "Light" lesson row (with its pdf link)
#  - keep the existing "Forest Our Life Line" and "Market" rows, minus stray data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean so no leftover cells/values remain from the old layout.
$ws.Cells.Clear()

$data = @(
    @("Grade", "Subjects", "Lesson", "Quiz", "Worksheet", "Flashcard"),
    @(1, "English", "Lesson 1", "", "", ""),
    @(1, "Math", "Addition", "", "", ""),
    @(2, "Science", "Plants", "", "", ""),
    @(7, "English", "Tenses", "", "", ""),
    @(7, "Math", "Ratio & Propotion", "", "", ""),
    @(7, "Science", "Light", "/Light.pdf", "", ""),
    @(7, "Science", "Forest Our Life Line", "https://quizizz.com/join?gc=08539312", "/forest.pdf", ""),
    @(7, "Social", "Market", "", "", "")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($val -ne "") {
            $ws.Cells.Item($r + 1, $c + 1).Value = $val
        }
    }
}

# Match the saved selection state from the edit.
$ws.Range("J11").Select()

Write-Host "done"
